$wb = $excel.ActiveWorkbook

# Add the new "Television" sheet at the end of the workbook (after the last
# existing sheet), so it becomes sheetId 3 / the 3rd tab and the active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Television"

# Populate the dialogue cells in the order they were originally typed so the
# shared-string table ends up in the same order as the target workbook.
$ws.Range("B5").Value = "FERN: Well, the trainers that are fighting - years ago they were huddled around a television, and watching, and wondering."
$ws.Range("B3").Value = "SCULLY: …and if you think about it, this year's Tournament is really a global one. Red, from Kanto. Silver, from Johto. Trainers from Unova, and Aloha, and Galar. "
$ws.Range("B8").Value = "SCULLY: On the scoreboard in Snowpoint Temple, it is `$CURRENTTIME in the city of miracles, Snowpoint City. And a crowd of 29,139 just sitting in to see the only trainer in tournament"
$ws.Range("B9").Value = "history to manuever four no-knock games, and she's done it four straight years. And now she's capped it: on her fourth no-knock victory, she made it a perfect match.`""
$ws.Range("B10").Value = "wow that doesn't work…sorry vin"
$ws.Range("B12").Value = "SCULLY: …her Typhlosion uses Burn Blast, lands a critical hit on Gengar, and she has done it! If you have a trainer hat, throw it to the sky!"
$ws.Range("B4").Value = "You think about the young kids from these regions, who are huddled around TVs with their families, watching this series. These trainers are inspiring the next wave of global talent."

# Match the final selection left on the new sheet.
$null = $ws.Range("B14").Select()
